$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# --- Add two rows to the "Consolidado_Tablas_Madre" table (extends table ref,
#     autoFilter and sheet dimension from L146 to L148) ---
$null = $lo.ListRows.Add()
$null = $lo.ListRows.Add()

# ===================== Row 147 =====================
# Seed formatting by copying the previous data row, then overwrite H/I with the
# exact cell styles used by the source rows (H147 -> "Hyperlink", I147 -> "Hipervínculo").
$ws.Range("A146:L146").Copy()
$ws.Range("A147:L147").PasteSpecial(-4122)
$ws.Cells.Item(147, 8).Style = "Hyperlink"
$ws.Cells.Item(147, 9).Style = "Hipervínculo"

$ws.Cells.Item(147, 1).Value = "20"
$ws.Cells.Item(147, 2).Value = "Pueblos_Comuna_Edad"
$ws.Cells.Item(147, 3).Value = ""
$ws.Cells.Item(147, 4).Value = "Persona"
$ws.Cells.Item(147, 5).Value = "Indefinido (Decenio)"
$ws.Cells.Item(147, 6).Value = "Manual"
$ws.Cells.Item(147, 7).Value = "Instituto Nacional de Estadísticas"
$ws.Cells.Item(147, 10).Value = "Andrés Sebastian"
$ws.Cells.Item(147, 11).Value = "Sociedad"
$ws.Cells.Item(147, 12).Value = "Demografía"

# ===================== Row 148 =====================
$ws.Range("A146:L146").Copy()
$ws.Range("A148:L148").PasteSpecial(-4122)
$ws.Cells.Item(148, 8).Style = "Hipervínculo"
$ws.Cells.Item(148, 9).Style = "Hipervínculo"

$ws.Cells.Item(148, 1).Value = "20.01"
$ws.Cells.Item(148, 2).Value = "Pueblos_Comuna_Genero"
$ws.Cells.Item(148, 3).Value = ""
$ws.Cells.Item(148, 4).Value = "Persona"
$ws.Cells.Item(148, 5).Value = "Indefinido (Decenio)"
$ws.Cells.Item(148, 6).Value = "Manual"
$ws.Cells.Item(148, 7).Value = "Instituto Nacional de Estadísticas"
$ws.Cells.Item(148, 10).Value = "Andrés Sebastian"
$ws.Cells.Item(148, 11).Value = "Sociedad"
$ws.Cells.Item(148, 12).Value = "Demografía"

# Move the active selection to A149 (just below the newly-added rows), matching
# the post-edit cursor position recorded in the workbook view.
$ws.Range("A149").Select()
